# Update cryptocurrency price/volume figures (and the EnergySwap/PancakeSwap row swap)
# to match the refreshed coinranking.com snapshot.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '24.241.70'
$ws.Range("E2").Value = '  +13.57%  '
$ws.Range("D3").Value = '1.674.06'
$ws.Range("E4").Value = '  -0.33%  '
$ws.Range("D5").Value = '''309.11'
$ws.Range("E5").Value = '  +9.35%  '
$ws.Range("D6").Value = '''0.9983'
$ws.Range("E6").Value = '  +3.03%  '
$ws.Range("D7").Value = '''0.3740'
$ws.Range("E7").Value = '  +3.12%  '
$ws.Range("D8").Value = '''0.3446'
$ws.Range("E8").Value = '  +7.69%  '
$ws.Range("D9").Value = '''47.69'
$ws.Range("E9").Value = '  +16.40%  '
$ws.Range("D10").Value = '''1.182'
$ws.Range("E10").Value = '  +5.90%  '
$ws.Range("D11").Value = '''0.07301'
$ws.Range("E11").Value = '  +5.41%  '
$ws.Range("D12").Value = '''0.9997'
$ws.Range("E12").Value = '  -0.16%  '
$ws.Range("D13").Value = '''20.46'
$ws.Range("E13").Value = '  +7.95%  '
$ws.Range("E14").Value = '  +6.70%  '
$ws.Range("D15").Value = '''6.768'
$ws.Range("E15").Value = '  +5.23%  '
$ws.Range("D16").Value = '1.677.46'
$ws.Range("E16").Value = '  +8.63%  '
$ws.Range("D17").Value = '''0.00001112'
$ws.Range("E18").Value = '  +3.09%  '
$ws.Range("D19").Value = '''0.06722'
$ws.Range("E19").Value = '  +9.30%  '
$ws.Range("D20").Value = '''81.86'
$ws.Range("E20").Value = '  +11.83%  '
$ws.Range("E21").Value = '  +7.70%  '
$ws.Range("D22").Value = '''6.146'
$ws.Range("E22").Value = '  +6.89%  '
$ws.Range("E23").Value = '  +5.08%  '
$ws.Range("D24").Value = '24.208.14'
$ws.Range("E24").Value = '  +13.35%  '
$ws.Range("D25").Value = '''2.416'
$ws.Range("E25").Value = '  +3.98%  '
$ws.Range("D26").Value = '''3.356'
$ws.Range("E26").Value = '  -9.35%  '
$ws.Range("D27").Value = '''2.668'
$ws.Range("E27").Value = '  +16.72%  '
$ws.Range("D28").Value = '''151.65'
$ws.Range("E28").Value = '  +2.53%  '
$ws.Range("D29").Value = '''19.51'
$ws.Range("E29").Value = '  +9.63%  '
$ws.Range("D30").Value = '1.861.49'
$ws.Range("E30").Value = '  +8.51%  '
$ws.Range("D31").Value = '''127.04'
$ws.Range("E31").Value = '  +7.11%  '
$ws.Range("D32").Value = '''6.424'
$ws.Range("E32").Value = '  +21.76%  '
$ws.Range("D33").Value = '''4.122'
$ws.Range("E33").Value = '  +1.85%  '
$ws.Range("D34").Value = '''0.9973'
$ws.Range("E34").Value = '  +14.26%  '
$ws.Range("D35").Value = '''1.768'
$ws.Range("E35").Value = '  +15.11%  '
$ws.Range("D36").Value = '''0.08521'
$ws.Range("E36").Value = '  +5.79%  '
$ws.Range("D37").Value = '''12.57'
$ws.Range("E37").Value = '  +16.26%  '
$ws.Range("D38").Value = '''0.06478'
$ws.Range("E38").Value = '  +10.32%  '
$ws.Range("E39").Value = '  +7.50%  '
$ws.Range("D40").Value = '''8.877'
$ws.Range("E40").Value = '  +10.69%  '
$ws.Range("D41").Value = '''0.02357'
$ws.Range("E41").Value = '  +10.65%  '
$ws.Range("D42").Value = '''1.282'
$ws.Range("E42").Value = '  +5.56%  '
$ws.Range("D43").Value = '''0.2154'
$ws.Range("D44").Value = '''0.6186'
$ws.Range("E44").Value = '  +11.95%  '
$ws.Range("D45").Value = '''0.9976'
$ws.Range("E45").Value = '  +2.96%  '
$ws.Range("B46").Value = 'EnergySwap'
$ws.Range("C46").Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range("D46").Value = '''13.27'
$ws.Range("E46").Value = '  +5.74%  '
$ws.Range("B47").Value = 'PancakeSwap'
$ws.Range("C47").Value = 'https://coinranking.com/coin/ncYFcP709+pancakeswap-cake'
$ws.Range("D47").Value = '''3.808'
$ws.Range("E47").Value = '  +6.26%  '
$ws.Range("D48").Value = '''0.5959'
$ws.Range("E48").Value = '  +8.14%  '
$ws.Range("D49").Value = '''127.27'
$ws.Range("E49").Value = '  +4.03%  '
$ws.Range("D50").Value = '''2.033'
$ws.Range("E50").Value = '  +7.92%  '
$ws.Range("D51").Value = '''0.07167'
$ws.Range("E51").Value = '  +8.33%  '
